# Auto-generated: update Kujata market-price snapshot values across sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 764.2679000000001
$ws.Range("J17").Value = 764.2679000000001
$ws.Range("L17").Value = 2292.8037
$ws.Range("N17").Value = -2628.8037
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 4000
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 4000
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -4968
$ws.Range("H58").Value = 552.3333
$ws.Range("I58").Value = 326
$ws.Range("K58").Value = 978
$ws.Range("M58").Value = -828
$ws.Range("H69").Value = 3980
$ws.Range("J69").Value = 3980
$ws.Range("L69").Value = 11940
$ws.Range("N69").Value = -13688
$ws.Range("H72").Value = 3980
$ws.Range("J72").Value = 3980
$ws.Range("L72").Value = 35820
$ws.Range("N72").Value = -44556
$ws.Range("H94").Value = 2580.8
$ws.Range("I94").Value = 2580.8
$ws.Range("K94").Value = 2580.8
$ws.Range("M94").Value = -2129.8
$ws.Range("H138").Value = 1524.6086
$ws.Range("I138").Value = 1095.8788
$ws.Range("J138").Value = 1764.4067
$ws.Range("K138").Value = 3287.6364
$ws.Range("L138").Value = 5293.2201
$ws.Range("M138").Value = 1852.3636
$ws.Range("N138").Value = -15573.2201

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2766.5334
$ws.Range("J88").Value = 2792.7856
$ws.Range("L88").Value = 2792.7856
$ws.Range("N88").Value = -3604.7856
$ws.Range("H91").Value = 2766.5334
$ws.Range("J91").Value = 2792.7856
$ws.Range("L91").Value = 2792.7856
$ws.Range("N91").Value = -5600.7856
$ws.Range("H97").Value = 205.95
$ws.Range("I97").Value = 211.57895
$ws.Range("J97").Value = 99
$ws.Range("K97").Value = 211.57895
$ws.Range("L97").Value = 99
$ws.Range("M97").Value = 284.42105
$ws.Range("N97").Value = -1091

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4722.722
$ws.Range("I86").Value = 5231.5386
$ws.Range("J86").Value = 3399.8
$ws.Range("K86").Value = 5231.5386
$ws.Range("L86").Value = 3399.8
$ws.Range("M86").Value = -4108.5386
$ws.Range("N86").Value = -5645.8
$ws.Range("H89").Value = 4722.722
$ws.Range("I89").Value = 5231.5386
$ws.Range("J89").Value = 3399.8
$ws.Range("K89").Value = 26157.693
$ws.Range("L89").Value = 16999
$ws.Range("M89").Value = -20541.693
$ws.Range("N89").Value = -28231
$ws.Range("H97").Value = 23575
$ws.Range("I97").Value = 9766.666999999999
$ws.Range("J97").Value = 65000
$ws.Range("K97").Value = 9766.666999999999
$ws.Range("L97").Value = 65000
$ws.Range("M97").Value = -8775.666999999999
$ws.Range("N97").Value = -66982

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 492.77777
$ws.Range("I22").Value = 362.14285
$ws.Range("J22").Value = 950
$ws.Range("K22").Value = 362.14285
$ws.Range("L22").Value = 950
$ws.Range("M22").Value = -12.14285000000001
$ws.Range("N22").Value = -1650
$ws.Range("H58").Value = 1660.3478
$ws.Range("I58").Value = 1377.1177
$ws.Range("K58").Value = 1377.1177
$ws.Range("M58").Value = -1174.1177
$ws.Range("H136").Value = 1660.3478
$ws.Range("I136").Value = 1377.1177
$ws.Range("K136").Value = 4131.3531
$ws.Range("M136").Value = -1581.3531

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 6054
$ws.Range("I63").Value = 4262
$ws.Range("K63").Value = 12786
$ws.Range("M63").Value = -12037
$ws.Range("H66").Value = 6054
$ws.Range("I66").Value = 4262
$ws.Range("K66").Value = 38358
$ws.Range("M66").Value = -34614
$ws.Range("H87").Value = 2240.75
$ws.Range("I87").Value = 1605.2
$ws.Range("J87").Value = 3300
$ws.Range("K87").Value = 4815.6
$ws.Range("L87").Value = 9900
$ws.Range("M87").Value = -3567.6
$ws.Range("N87").Value = -12396
$ws.Range("H90").Value = 2240.75
$ws.Range("I90").Value = 1605.2
$ws.Range("J90").Value = 3300
$ws.Range("K90").Value = 14446.8
$ws.Range("L90").Value = 29700
$ws.Range("M90").Value = -8206.800000000001
$ws.Range("N90").Value = -42180
$ws.Range("H113").Value = 685.08
$ws.Range("J113").Value = 685.08
$ws.Range("L113").Value = 2055.24
$ws.Range("N113").Value = -6395.24
$ws.Range("H114").Value = 511.35715
$ws.Range("I114").Value = 569.2
$ws.Range("J114").Value = 479.22223
$ws.Range("K114").Value = 1707.6
$ws.Range("L114").Value = 1437.66669
$ws.Range("M114").Value = 1546.4
$ws.Range("N114").Value = -7945.66669
$ws.Range("H129").Value = 34722904
$ws.Range("I129").Value = 37037500
$ws.Range("K129").Value = 111112500
$ws.Range("M129").Value = -111107500
$ws.Range("H131").Value = 14707098
$ws.Range("J131").Value = 1244.7122
$ws.Range("L131").Value = 3734.1366
$ws.Range("N131").Value = -13814.1366

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 2000
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()
$ws.Range("H80").Value = 3635.625
$ws.Range("I80").Value = 2580.8333
$ws.Range("K80").Value = 2580.8333
$ws.Range("M80").Value = -1582.8333
$ws.Range("H83").Value = 3635.625
$ws.Range("I83").Value = 2580.8333
$ws.Range("K83").Value = 12904.1665
$ws.Range("M83").Value = -7912.166499999999
$ws.Range("H126").Value = 2801.375
$ws.Range("I126").Value = 3102.75
$ws.Range("J126").Value = 2500
$ws.Range("K126").Value = 9308.25
$ws.Range("L126").Value = 7500
$ws.Range("M126").Value = -6838.25
$ws.Range("N126").Value = -12440
$ws.Range("H132").Value = 2669.5557
$ws.Range("I132").Value = 2243.8333
$ws.Range("J132").Value = 3521
$ws.Range("K132").Value = 6731.499899999999
$ws.Range("L132").Value = 10563
$ws.Range("M132").Value = -4201.499899999999
$ws.Range("N132").Value = -15623

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4500.3335
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 4500.3335
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 4500.3335
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -5090.3335
$ws.Range("H27").Value = 4500.3335
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 4500.3335
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 4500.3335
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -4714.3335

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 700005
$ws.Range("I14").Value = 700005
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 700005
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -699837
$ws.Range("N14").ClearContents()
$ws.Range("H81").Value = 655.25
$ws.Range("I81").Value = 641.6667
$ws.Range("J81").Value = 696
$ws.Range("K81").Value = 1283.3334
$ws.Range("L81").Value = 1392
$ws.Range("M81").Value = -222.3334
$ws.Range("N81").Value = -3514
$ws.Range("H84").Value = 655.25
$ws.Range("I84").Value = 641.6667
$ws.Range("J84").Value = 696
$ws.Range("K84").Value = 6416.666999999999
$ws.Range("L84").Value = 6960
$ws.Range("M84").Value = -1112.666999999999
$ws.Range("N84").Value = -17568

Write-Host "Kujata profit snapshot updated"